$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(" <sup>imprimer</sup>", $true, $false, $false, $false, $false,
              $true, 1, $false, "", 2)
